$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D stores crypto prices as plain text (many values, e.g. "29.144.26",
# use dots as thousands separators and are not valid numbers). For the handful
# of updated prices that do look like ordinary decimals (e.g. "0.9991"), force
# the Text number format on just those cells first so Excel keeps them as text
# instead of auto-converting them to numbers, matching the rest of the column.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D17", "D20", "D21", "D23", "D25", "D26", "D30", "D31", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.144.26"
$ws.Range("D3").Value = "1.833.31"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "241.54"
$ws.Range("D6").Value = "0.6583"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D8").Value = "0.07413"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.2927"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").Value = "22.89"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.840.64"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "4.993"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "0.6674"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").Value = "6.115"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "0.000008614"
$ws.Range("E17").Value = "  +4.68%  "
$ws.Range("D18").Value = "29.143.44"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "2.084.27"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "226.59"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "7.113"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "161.48"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "0.1405"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "4.102"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("D31").Value = "4.048"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "0.05269"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").Value = "1.866"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").Value = "1.143"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "2.654"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "1.305.71"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "2.741"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "0.9158"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "6.066"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "102.22"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D45").Value = "0.08119"
$ws.Range("E45").Value = "  +10.86%  "
$ws.Range("D46").Value = "1.982.68"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "0.5129"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "63.70"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "1.749"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "0.05841"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "6.769"
$ws.Range("E51").Value = "  -0.95%  "
